$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1088.4131
$ws.Range("I17").Value = 1150
$ws.Range("J17").Value = 1085.6136
$ws.Range("K17").Value = 3450
$ws.Range("L17").Value = 3256.8408
$ws.Range("M17").Value = -3282
$ws.Range("N17").Value = -3592.8408
$ws.Range("H32").Value = 4821.077
$ws.Range("I32").Value = 2500
$ws.Range("J32").Value = 4913.92
$ws.Range("K32").Value = 2500
$ws.Range("L32").Value = 4913.92
$ws.Range("M32").Value = -2174
$ws.Range("N32").Value = -5565.92
$ws.Range("H98").Value = 1254.579
$ws.Range("I98").Value = 1254.579
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1254.579
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 243.421
$ws.Range("H99").Value = 503.66666
$ws.Range("I99").Value = 323.75
$ws.Range("J99").Value = 863.5
$ws.Range("K99").Value = 971.25
$ws.Range("L99").Value = 2590.5
$ws.Range("M99").Value = 526.75
$ws.Range("N99").Value = -5586.5
$ws.Range("H122").Value = 1254.579
$ws.Range("I122").Value = 1254.579
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3763.737
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1313.737
$ws.Range("H137").Value = 73742.75999999999
$ws.Range("I137").Value = 150248
$ws.Range("J137").Value = 3122.5386
$ws.Range("K137").Value = 450744
$ws.Range("L137").Value = 9367.6158
$ws.Range("M137").Value = -448194
$ws.Range("N137").Value = -14467.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 5186.75
$ws.Range("I19").Value = 3213.4285
$ws.Range("J19").Value = 19000
$ws.Range("K19").Value = 3213.4285
$ws.Range("L19").Value = 19000
$ws.Range("M19").Value = -2984.4285
$ws.Range("N19").Value = -19458
$ws.Range("H32").Value = 4656.108
$ws.Range("I32").Value = 3069.258
$ws.Range("J32").Value = 12854.833
$ws.Range("K32").Value = 3069.258
$ws.Range("L32").Value = 12854.833
$ws.Range("M32").Value = -2782.258
$ws.Range("N32").Value = -13428.833
$ws.Range("H61").Value = 1569.8125
$ws.Range("I61").Value = 1276.8182
$ws.Range("J61").Value = 2214.4
$ws.Range("K61").Value = 1276.8182
$ws.Range("L61").Value = 2214.4
$ws.Range("M61").Value = -1064.8182
$ws.Range("N61").Value = -2638.4
$ws.Range("H74").Value = 73522.94500000001
$ws.Range("I74").Value = 6990.8438
$ws.Range("J74").Value = 377669.72
$ws.Range("K74").Value = 6990.8438
$ws.Range("L74").Value = 377669.72
$ws.Range("M74").Value = -6116.8438
$ws.Range("N74").Value = -379417.72
$ws.Range("H77").Value = 73522.94500000001
$ws.Range("I77").Value = 6990.8438
$ws.Range("J77").Value = 377669.72
$ws.Range("K77").Value = 34954.219
$ws.Range("L77").Value = 1888348.6
$ws.Range("M77").Value = -30586.219
$ws.Range("N77").Value = -1897084.6
$ws.Range("H88").Value = 1454.8182
$ws.Range("I88").Value = 475.5
$ws.Range("J88").Value = 2014.4286
$ws.Range("K88").Value = 475.5
$ws.Range("L88").Value = 2014.4286
$ws.Range("M88").Value = -69.5
$ws.Range("N88").Value = -2826.4286
$ws.Range("H91").Value = 1454.8182
$ws.Range("I91").Value = 475.5
$ws.Range("J91").Value = 2014.4286
$ws.Range("K91").Value = 475.5
$ws.Range("L91").Value = 2014.4286
$ws.Range("M91").Value = 928.5
$ws.Range("N91").Value = -4822.4286
$ws.Range("H132").Value = 2273.6667
$ws.Range("I132").Value = 1853.8148
$ws.Range("J132").Value = 4163
$ws.Range("K132").Value = 5561.4444
$ws.Range("L132").Value = 12489
$ws.Range("M132").Value = -3031.4444
$ws.Range("N132").Value = -17549
$ws.Range("H136").Value = 1569.8125
$ws.Range("I136").Value = 1276.8182
$ws.Range("J136").Value = 2214.4
$ws.Range("K136").Value = 3830.4546
$ws.Range("L136").Value = 6643.200000000001
$ws.Range("M136").Value = -1280.4546
$ws.Range("N136").Value = -11743.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7153261.5
$ws.Range("I86").Value = 9103242
$ws.Range("J86").Value = 3333.3333
$ws.Range("K86").Value = 9103242
$ws.Range("L86").Value = 3333.3333
$ws.Range("M86").Value = -9102119
$ws.Range("N86").Value = -5579.3333
$ws.Range("H89").Value = 7153261.5
$ws.Range("I89").Value = 9103242
$ws.Range("J89").Value = 3333.3333
$ws.Range("K89").Value = 45516210
$ws.Range("L89").Value = 16666.6665
$ws.Range("M89").Value = -45510594
$ws.Range("N89").Value = -27898.6665
$ws.Range("H134").Value = 2656.9592
$ws.Range("I134").Value = 1034.8889
$ws.Range("J134").Value = 7148.846
$ws.Range("K134").Value = 3104.6667
$ws.Range("L134").Value = 21446.538
$ws.Range("M134").Value = -569.6666999999998
$ws.Range("N134").Value = -26516.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31350.773
$ws.Range("I31").Value = 1352.5294
$ws.Range("J31").Value = 67777.21000000001
$ws.Range("K31").Value = 1352.5294
$ws.Range("L31").Value = 67777.21000000001
$ws.Range("M31").Value = -1057.5294
$ws.Range("N31").Value = -68367.21000000001
$ws.Range("H34").Value = 31350.773
$ws.Range("I34").Value = 1352.5294
$ws.Range("J34").Value = 67777.21000000001
$ws.Range("K34").Value = 1352.5294
$ws.Range("L34").Value = 67777.21000000001
$ws.Range("M34").Value = -1150.5294
$ws.Range("N34").Value = -68181.21000000001
$ws.Range("H58").Value = 1525.4474
$ws.Range("I58").Value = 1171.8334
$ws.Range("J58").Value = 2851.5
$ws.Range("K58").Value = 1171.8334
$ws.Range("L58").Value = 2851.5
$ws.Range("M58").Value = -968.8334
$ws.Range("N58").Value = -3257.5
$ws.Range("H108").Value = 35093.875
$ws.Range("I108").Value = 23999.5
$ws.Range("J108").Value = 38792
$ws.Range("K108").Value = 23999.5
$ws.Range("L108").Value = 38792
$ws.Range("M108").Value = -20159.5
$ws.Range("N108").Value = -46472
$ws.Range("H136").Value = 1525.4474
$ws.Range("I136").Value = 1171.8334
$ws.Range("J136").Value = 2851.5
$ws.Range("K136").Value = 3515.5002
$ws.Range("L136").Value = 8554.5
$ws.Range("M136").Value = -965.5001999999999
$ws.Range("N136").Value = -13654.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 84713.836
$ws.Range("I5").Value = 517.6
$ws.Range("J5").Value = 144854
$ws.Range("K5").Value = 1552.8
$ws.Range("L5").Value = 434562
$ws.Range("M5").Value = -1440.8
$ws.Range("N5").Value = -434786
$ws.Range("H56").Value = 10422533
$ws.Range("I56").Value = 10422533
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 10422533
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -10422003
$ws.Range("H74").Value = 4931.5
$ws.Range("I74").Value = 1575.3334
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 4726.0002
$ws.Range("L74").Value = 45000
$ws.Range("M74").Value = -3665.0002
$ws.Range("N74").Value = -47122
$ws.Range("H77").Value = 4931.5
$ws.Range("I77").Value = 1575.3334
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 14178.0006
$ws.Range("L77").Value = 135000
$ws.Range("M77").Value = -8874.000599999999
$ws.Range("N77").Value = -145608
$ws.Range("H107").Value = 235.22728
$ws.Range("I107").Value = 172.11765
$ws.Range("J107").Value = 449.8
$ws.Range("K107").Value = 516.35295
$ws.Range("L107").Value = 1349.4
$ws.Range("M107").Value = 1403.64705
$ws.Range("N107").Value = -5189.4
$ws.Range("H122").Value = 953.7
$ws.Range("I122").Value = 1133.3334
$ws.Range("J122").Value = 876.7143
$ws.Range("K122").Value = 10200.0006
$ws.Range("L122").Value = 7890.428699999999
$ws.Range("M122").Value = -7750.000599999999
$ws.Range("N122").Value = -12790.4287
$ws.Range("H135").Value = 84713.836
$ws.Range("I135").Value = 517.6
$ws.Range("J135").Value = 144854
$ws.Range("K135").Value = 4658.400000000001
$ws.Range("L135").Value = 1303686
$ws.Range("M135").Value = -2123.400000000001
$ws.Range("N135").Value = -1308756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 17833.334
$ws.Range("I58").Value = 9250
$ws.Range("J58").Value = 35000
$ws.Range("K58").Value = 9250
$ws.Range("L58").Value = 35000
$ws.Range("M58").Value = -8973
$ws.Range("N58").Value = -35554
$ws.Range("H132").Value = 2866.513
$ws.Range("I132").Value = 2400.9092
$ws.Range("J132").Value = 5427.3335
$ws.Range("K132").Value = 7202.7276
$ws.Range("L132").Value = 16282.0005
$ws.Range("M132").Value = -4672.7276
$ws.Range("N132").Value = -21342.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4612.314
$ws.Range("I132").Value = 3632.75
$ws.Range("J132").Value = 6963.2666
$ws.Range("K132").Value = 10898.25
$ws.Range("L132").Value = 20889.7998
$ws.Range("M132").Value = -8368.25
$ws.Range("N132").Value = -25949.7998
$ws.Range("H136").Value = 65004.453
$ws.Range("I136").Value = 87010.28999999999
$ws.Range("J136").Value = 6322.222
$ws.Range("K136").Value = 261030.87
$ws.Range("L136").Value = 18966.666
$ws.Range("M136").Value = -258480.87
$ws.Range("N136").Value = -24066.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6552.754
$ws.Range("I62").Value = 3249.9644
$ws.Range("J62").Value = 9355.120999999999
$ws.Range("K62").Value = 3249.9644
$ws.Range("L62").Value = 9355.120999999999
$ws.Range("M62").Value = -2625.9644
$ws.Range("N62").Value = -10603.121
$ws.Range("H65").Value = 6552.754
$ws.Range("I65").Value = 3249.9644
$ws.Range("J65").Value = 9355.120999999999
$ws.Range("K65").Value = 16249.822
$ws.Range("L65").Value = 46775.605
$ws.Range("M65").Value = -13129.822
$ws.Range("N65").Value = -53015.605
$ws.Range("H120").Value = 40999
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 40999
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 40999
$ws.Range("N120").Value = -50675
$ws.Range("H129").Value = 39349.5
$ws.Range("I129").Value = 39000
$ws.Range("J129").Value = 39699
$ws.Range("K129").Value = 39000
$ws.Range("L129").Value = 39699
$ws.Range("M129").Value = -34000
$ws.Range("N129").Value = -49699
$ws.Range("H132").Value = 21977186
$ws.Range("I132").Value = 23810562
$ws.Range("J132").Value = 2726722.8
$ws.Range("K132").Value = 71431686
$ws.Range("L132").Value = 8180168.399999999
$ws.Range("M132").Value = -71429156
$ws.Range("N132").Value = -8185228.399999999
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 1528.0555
$ws.Range("J136").Value = 4831.6665
$ws.Range("K136").Value = 4584.166499999999
$ws.Range("L136").Value = 14494.9995
$ws.Range("M136").Value = -2034.166499999999
$ws.Range("N136").Value = -19594.9995
